$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.807.24"
$ws.Cells.Item(2, 5).Value = "  -2.00%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.800.28"
$ws.Cells.Item(3, 5).Value = "  -1.44%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "308.65"
$ws.Cells.Item(5, 5).Value = "  -1.80%  "
$ws.Cells.Item(6, 5).Value = "  +0.05%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4640"
$ws.Cells.Item(7, 5).Value = "  +3.78%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3681"
$ws.Cells.Item(8, 5).Value = "  -2.03%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07342"
$ws.Cells.Item(9, 5).Value = "  -1.74%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8654"
$ws.Cells.Item(10, 5).Value = "  -2.43%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.35"
$ws.Cells.Item(11, 5).Value = "  -3.16%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.894.58"
$ws.Cells.Item(12, 5).Value = "  +3.73%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.344"
$ws.Cells.Item(13, 5).Value = "  -1.57%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.512"
$ws.Cells.Item(14, 5).Value = "  -3.60%  "
$ws.Cells.Item(15, 5).Value = "  -1.40%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "91.22"
$ws.Cells.Item(16, 5).Value = "  -2.69%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.002"
$ws.Cells.Item(17, 5).Value = "  +0.13%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008685"
$ws.Cells.Item(18, 5).Value = "  -1.12%  "
$ws.Cells.Item(19, 5).Value = "  +0.08%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.60"
$ws.Cells.Item(20, 5).Value = "  -3.68%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "26.815.74"
$ws.Cells.Item(21, 5).Value = "  -1.96%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.287"
$ws.Cells.Item(22, 5).Value = "  -2.31%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.58"
$ws.Cells.Item(23, 5).Value = "  -3.43%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.092.10"
$ws.Cells.Item(24, 5).Value = "  +1.65%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.904"
$ws.Cells.Item(25, 5).Value = "  -3.32%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "151.42"
$ws.Cells.Item(26, 5).Value = "  +0.06%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.31"
$ws.Cells.Item(27, 5).Value = "  -2.39%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.120"
$ws.Cells.Item(28, 5).Value = "  -8.25%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.226"
$ws.Cells.Item(29, 5).Value = "  -3.13%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "115.74"
$ws.Cells.Item(30, 5).Value = "  -1.78%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.08902"
$ws.Cells.Item(31, 5).Value = "  +0.22%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.7541"
$ws.Cells.Item(32, 5).Value = "  -4.53%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.932"
$ws.Cells.Item(33, 5).Value = "  +0.14%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.144"
$ws.Cells.Item(34, 5).Value = "  -5.20%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.439"
$ws.Cells.Item(35, 5).Value = "  -3.59%  "
$ws.Cells.Item(36, 5).Value = "  +0.07%  "
$ws.Cells.Item(37, 5).Value = "  -0.80%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01947"
$ws.Cells.Item(38, 5).Value = "  -2.24%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.05237"
$ws.Cells.Item(39, 5).Value = "  -1.35%  "
$ws.Cells.Item(40, 5).Value = "  +2.56%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "7.178"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.5262"
$ws.Cells.Item(42, 5).Value = "  -1.98%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.342"
$ws.Cells.Item(43, 5).Value = "  +0.83%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.1656"
$ws.Cells.Item(44, 5).Value = "  -3.68%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "8.447"
$ws.Cells.Item(45, 5).Value = "  -2.58%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4989"
$ws.Cells.Item(46, 5).Value = "  -2.81%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "10.30"
$ws.Cells.Item(47, 5).Value = "  -3.01%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "103.97"
$ws.Cells.Item(48, 5).Value = "  -1.17%  "
$ws.Cells.Item(49, 5).Value = "  +0.03%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.659"
$ws.Cells.Item(50, 5).Value = "  -2.30%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06284"
$ws.Cells.Item(51, 5).Value = "  -1.99%  "
